$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(3).Delete()
$ws.Rows(3).Delete()
